# "display parallel data for tree"
#
# 1. On the existing "Memory Usage" sheet: replace the old (sequential)
#    memory-usage data in A2:F6 with the new parallel-run data, now laid
#    out in columns G:J (rows 2-4).
# 2. Add a new "Execution Times" sheet with the matching (s) headers and
#    one row of timing data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Memory Usage sheet: clear old data block and write the new values ---
$ws1.Range("A2:F6").ClearContents()

$ws1.Range("G2").Value = 3067096.0
$ws1.Range("H2").Value = 3697080.0
$ws1.Range("I2").Value = 1331752.0
$ws1.Range("J2").Value = 1331736.0

$ws1.Range("G3").Value = 1285760.0
$ws1.Range("H3").Value = 3833584.0
$ws1.Range("I3").Value = 2663464.0
$ws1.Range("J3").Value = 1331736.0

$ws1.Range("I4").Value = 0.0
$ws1.Range("J4").Value = 1331736.0

# Update the visible selection on the sheet to match the new data block.
$ws1.Range("A2:F6").Select() | Out-Null

# --- Add the "Execution Times" sheet after "Memory Usage" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Execution Times"

$headers = @(
  "Graph BFS 10000 (s)",
  "Graph DFS 10000 (s)",
  "Graph BFS 1000 (s)",
  "Graph DFS 1000 (s)",
  "Graph BFS 50000 (s)",
  "Graph DFS 50000 (s)",
  "Tree BFS 10000 (s)",
  "Tree DFS 10000 (s)",
  "Tree BFS 1000 (s)",
  "Tree DFS 1000 (s)",
  "Tree BFS 50000 (s)",
  "Tree DFS 50000 (s)"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws2.Range("A2").Value = 2663432
$ws2.Range("B2").Value = 2663432
$ws2.Range("C2").Value = 0.0
$ws2.Range("D2").Value = 0.0

# Leave "Memory Usage" as the active/selected sheet, as before.
$ws1.Activate() | Out-Null
